# Rename 'programs' sheet to 'time_variants' and update related selection state.

$wb = $excel.ActiveWorkbook

# 1. Rename the 'programs' worksheet to 'time_variants'.
$wsPrograms = $wb.Worksheets.Item("programs")
$wsPrograms.Name = "time_variants"

# 2. Change the dropdown value on model_attributes!B22 from 'scipy' to 'explicit'.
$wsAttrs = $wb.Worksheets.Item("model_attributes")
$wsAttrs.Range("B22").Value = "explicit"

# 3. Update the active sheet / selection state.
#    The active sheet moves from 'time_variants' (previously 'programs') back to
#    'model_attributes', with a new selection on model_attributes and a new
#    selection on time_variants (for when it is revisited).

$wsTimeVariants = $wb.Worksheets.Item("time_variants")
$wsTimeVariants.Range("K10").Select() | Out-Null

$wsAttrs.Select()
$wsAttrs.Range("C19").Select() | Out-Null
